$wb = $excel.ActiveWorkbook

# "Import Priorities" is the first sheet (sheetId=2, r:id=rId1)
$wsImportPriorities = $wb.Worksheets.Item(1)
# "Coupling Parameters" is the second sheet (sheetId=1, r:id=rId2)
$wsCouplingParameters = $wb.Worksheets.Item(2)

# Bump the import priority values (column B, rows 3-10) up by one.
$wsImportPriorities.Range("B3").Value = 9
$wsImportPriorities.Range("B4").Value = 8
$wsImportPriorities.Range("B5").Value = 7
$wsImportPriorities.Range("B6").Value = 6
$wsImportPriorities.Range("B7").Value = 5
$wsImportPriorities.Range("B8").Value = 4
$wsImportPriorities.Range("B9").Value = 3
$wsImportPriorities.Range("B10").Value = 2

# Switch the active/selected sheet from "Coupling Parameters" back to
# "Import Priorities", and move the selection on that sheet to C6.
$wsImportPriorities.Activate()
$wsImportPriorities.Range("C6").Select()

# Restore the workbook window's on-screen position/size (maximized layout).
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -16320
$win.Width = 29040
$win.Height = 15840
